$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 423.16666
$ws.Range("I135").Value = 370.72726
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 3336.54534
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -801.5453400000001
$ws.Range("N135").Value = -14070

$ws.Range("H139").Value = 25156
$ws.Range("J139").Value = 25156
$ws.Range("L139").Value = 25156
$ws.Range("N139").Value = -35436


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 61260.234
$ws.Range("I99").Value = 92537.17999999999
$ws.Range("K99").Value = 92537.17999999999
$ws.Range("M99").Value = -91039.17999999999

$ws.Range("H126").Value = 61260.234
$ws.Range("I126").Value = 92537.17999999999
$ws.Range("K126").Value = 277611.54
$ws.Range("M126").Value = -275141.54


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1000
$ws.Range("J13").Value = 1000
$ws.Range("L13").Value = 3000
$ws.Range("N13").Value = -3336

$ws.Range("H62").Value = 3550
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3550
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 10650
$ws.Range("N62").Value = -12022
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 3550
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3550
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 31950
$ws.Range("N65").Value = -38814
$ws.Range("M65").ClearContents()

$ws.Range("H69").Value = 2465.8333
$ws.Range("I69").Value = 991.4286
$ws.Range("J69").Value = 3072.9412
$ws.Range("K69").Value = 2974.2858
$ws.Range("L69").Value = 9218.8236
$ws.Range("M69").Value = -2163.2858
$ws.Range("N69").Value = -10840.8236

$ws.Range("H70").Value = 2640.2856
$ws.Range("I70").Value = 2096.4
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 6289.200000000001
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -5974.200000000001
$ws.Range("N70").Value = -12630

$ws.Range("H72").Value = 2465.8333
$ws.Range("I72").Value = 991.4286
$ws.Range("J72").Value = 3072.9412
$ws.Range("K72").Value = 8922.857399999999
$ws.Range("L72").Value = 27656.4708
$ws.Range("M72").Value = -4866.857399999999
$ws.Range("N72").Value = -35768.4708

$ws.Range("H73").Value = 2640.2856
$ws.Range("I73").Value = 2096.4
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 6289.200000000001
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -5197.200000000001
$ws.Range("N73").Value = -14184

$ws.Range("H75").Value = 1900
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 1975
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 5925
$ws.Range("M75").Value = -2002
$ws.Range("N75").Value = -7921

$ws.Range("H76").Value = 613
$ws.Range("I76").Value = 613
$ws.Range("K76").Value = 1839
$ws.Range("M76").Value = -1456

$ws.Range("H78").Value = 1900
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 1975
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 17775
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -27759

$ws.Range("H79").Value = 613
$ws.Range("I79").Value = 613
$ws.Range("K79").Value = 1839
$ws.Range("M79").Value = -513

$ws.Range("H82").Value = 7166.6665
$ws.Range("J82").Value = 7166.6665
$ws.Range("L82").Value = 21499.9995
$ws.Range("N82").Value = -22311.9995

$ws.Range("H85").Value = 7166.6665
$ws.Range("J85").Value = 7166.6665
$ws.Range("L85").Value = 21499.9995
$ws.Range("N85").Value = -24307.9995

$ws.Range("H88").Value = 5600
$ws.Range("J88").Value = 5600
$ws.Range("L88").Value = 16800
$ws.Range("N88").Value = -17656

$ws.Range("H91").Value = 5600
$ws.Range("J91").Value = 5600
$ws.Range("L91").Value = 16800
$ws.Range("N91").Value = -19764

$ws.Range("H94").Value = 7521.7144
$ws.Range("I94").Value = 3024
$ws.Range("J94").Value = 7867.6924
$ws.Range("K94").Value = 9072
$ws.Range("L94").Value = 23603.0772
$ws.Range("M94").Value = -8396
$ws.Range("N94").Value = -24955.0772

$ws.Range("H100").Value = 8432.5
$ws.Range("J100").Value = 8922.223
$ws.Range("L100").Value = 26766.669
$ws.Range("N100").Value = -28388.669

$ws.Range("H103").Value = 990.6923
$ws.Range("I103").Value = 255.1
$ws.Range("J103").Value = 3442.6667
$ws.Range("K103").Value = 765.3
$ws.Range("L103").Value = 10328.0001
$ws.Range("M103").Value = 113.7
$ws.Range("N103").Value = -12086.0001

$ws.Range("H112").Value = 2511
$ws.Range("I112").Value = 1639.6
$ws.Range("J112").Value = 3963.3333
$ws.Range("K112").Value = 4918.799999999999
$ws.Range("L112").Value = 11889.9999
$ws.Range("M112").Value = -3810.799999999999
$ws.Range("N112").Value = -14105.9999

$ws.Range("H116").Value = 5059.75
$ws.Range("I116").Value = 785.3333
$ws.Range("J116").Value = 6046.154
$ws.Range("K116").Value = 2355.9999
$ws.Range("L116").Value = 18138.462
$ws.Range("M116").Value = 1086.0001
$ws.Range("N116").Value = -25022.462

$ws.Range("H122").Value = 17242076
$ws.Range("I122").Value = 27778058
$ws.Range("J122").Value = 1379
$ws.Range("K122").Value = 250002522
$ws.Range("L122").Value = 12411
$ws.Range("M122").Value = -250000072
$ws.Range("N122").Value = -17311

$ws.Range("H123").Value = 2950
$ws.Range("I123").Value = 850
$ws.Range("J123").Value = 4000
$ws.Range("K123").Value = 2550
$ws.Range("L123").Value = 12000
$ws.Range("M123").Value = -100
$ws.Range("N123").Value = -16900

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3587.4167
$ws.Range("J7").Value = 3321.2856
$ws.Range("L7").Value = 3321.2856
$ws.Range("N7").Value = -3545.2856

$ws.Range("H40").Value = 3838.818
$ws.Range("I40").Value = 3805.5
$ws.Range("J40").Value = 3857.8572
$ws.Range("K40").Value = 3805.5
$ws.Range("L40").Value = 3857.8572
$ws.Range("M40").Value = -3669.5
$ws.Range("N40").Value = -4129.8572

$ws.Range("H126").Value = 3587.4167
$ws.Range("J126").Value = 3321.2856
$ws.Range("L126").Value = 9963.856800000001
$ws.Range("N126").Value = -14903.8568

$ws.Range("H127").Value = 44000
$ws.Range("J127").Value = 44000
$ws.Range("L127").Value = 44000
$ws.Range("N127").Value = -53920


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1010.871
$ws.Range("I126").Value = 508.6
$ws.Range("J126").Value = 1924.091
$ws.Range("K126").Value = 1525.8
$ws.Range("L126").Value = 5772.272999999999
$ws.Range("M126").Value = 944.1999999999998
$ws.Range("N126").Value = -10712.273
